$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1) Split the "href " run (paragraph 4) so the word "href" is wrapped
#        with spell-check proofErr markers, leaving the trailing space in
#        its own run. ---
$p4 = $d.Paragraphs(4)
$body4 = '<w:body>' `
  + '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">&lt;a </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>href</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>&#8211; hyper link (it will take you to another web page)</w:t></w:r>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' `
  + '</w:p></w:body>'
$p4.Range.InsertXML($pkgHeader + $body4 + $pkgFooter)

# --- 2) Split the "<img ..." paragraph, split the "Ol ..." paragraph, and
#        add the two new "TD" / "Tr" timetable paragraphs right after. ---
$p5 = $d.Paragraphs(5)
$p6 = $d.Paragraphs(6)
$combined = $d.Range($p5.Range.Start, $p6.Range.End)
$body56 = '<w:body>' `
  + '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>&lt;</w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>img</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> &#8211; image (insert an image)</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Ol</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> &gt; ordered list</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>TD &gt; stands for table Data, creating cells in rows</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Tr &gt; </w:t></w:r>' `
  + '</w:p>' `
  + '</w:body>'
$combined.InsertXML($pkgHeader + $body56 + $pkgFooter)
